# Reorder the comma-separated "Recorded By" names in column G so the
# first-listed recorder is moved to the end of the list (left-rotation
# of the comma-separated values). Cells with only a single value (no
# comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "*,*") {
        $parts = $val -split ",\s*"

        $rotated = @()
        for ($i = 1; $i -lt $parts.Count; $i++) {
            $rotated += $parts[$i]
        }
        $rotated += $parts[0]

        $cell.Value = [string]::Join(", ", $rotated)
    }
}
